# Updates the cryptos price/volume table with refreshed values.
# Note: several "Price" (column D) values look like plain numbers
# (e.g. "515.95", "6.00", "0.999"). A leading apostrophe is used for
# those so Excel stores them as text (matching the workbook's existing
# inline/shared-string cell type) instead of silently converting them
# to floating point numbers and losing the original formatting
# (trailing zeros, thousand-separator-looking dots, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.903.83'
$ws.Range('D3').Value = '3.061.53'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''515.95'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = '''141.51'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').Value = '''7.29'
$ws.Range('E9').Value = '  +2.22%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').Value = '3.584.24'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').Value = '''26.33'
$ws.Range('E14').Value = '  +3.81%  '
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '57.917.87'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').Value = '3.057.15'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '''6.11'
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('E19').Value = '  -2.58%  '
$ws.Range('D20').Value = '''8.16'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').Value = '''330.76'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '0.0₃0905'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = '''6.46'
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('E29').Value = '  +5.36%  '
$ws.Range('D30').Value = '''1.81'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').Value = '''1.20'
$ws.Range('E31').Value = '  +3.30%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').Value = '''154.91'
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('B35').Value = 'EnergySwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D35').Value = '''27.30'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '''6.00'
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('E37').Value = '  +3.27%  '
$ws.Range('D38').Value = '''0.0678'
$ws.Range('E38').Value = '  +2.12%  '
$ws.Range('D39').Value = '3.102.82'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('D41').Value = '''36.67'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').Value = '''0.656'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '2.301.27'
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('E45').Value = '  +4.19%  '
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').Value = '''20.82'
$ws.Range('E47').Value = '  +4.95%  '
$ws.Range('D48').Value = '''0.939'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '''0.733'
$ws.Range('E50').Value = '  +8.57%  '
$ws.Range('D51').Value = '''254.35'
$ws.Range('E51').Value = '  +9.53%  '
